$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.622.92"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "1.818.70"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.582"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.23%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "35.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.301"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0697"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").Value = "2.083.52"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").Value = "1.818.79"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.648"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("D16").Value = "34.600.01"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.98%  "
$ws.Range("D19").Value = "0.0₃0801"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "246.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.27%  "
$ws.Range("E22").Value = "  +0.12%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "173.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.54"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0531"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.82%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.84"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("D35").Value = "1.402.51"
$ws.Range("E35").Value = "  -1.87%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.56"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.30%  "
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.679"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0191"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.86"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "83.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.953"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.82%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.93"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.56%  "
$ws.Range("E45").Value = "  +2.60%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0517"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.23%  "
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").Value = "1.982.98"
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.02%  "
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("E51").Value = "  +0.15%  "
